$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.406.16"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.801.52"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'225.23"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "'0.598"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "'36.03"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("D9").Value = "'0.290"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "2.062.25"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'11.21"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "1.799.49"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "34.369.46"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'4.42"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "'68.58"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'242.43"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "0.0₃0774"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "'11.31"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +5.53%  "
$ws.Range("D25").Value = "'171.09"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'7.84"
$ws.Range("E26").Value = "  +4.72%  "
$ws.Range("D27").Value = "'17.38"
$ws.Range("E27").Value = "  +4.43%  "
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.23"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.79"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").Value = "1.363.83"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'2.36"
$ws.Range("E38").Value = "  -7.63%  "
$ws.Range("D39").Value = "'0.0186"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "'81.16"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E44").Value = "  +4.91%  "
$ws.Range("D45").Value = "'13.22"
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("D46").Value = "'0.0501"
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").Value = "1.964.02"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "'101.74"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").Value = "0.0₆0123"
$ws.Range("E51").Value = "  -0.10%  "
